$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $text) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.ClearFormats()
}

# --- Remove column C entirely (duplicated data no longer needed) ---
$ws.Columns.Item(3).Delete()

# --- Update existing B-column values (rows 5-14) ---
$ws.Range("B5").Value = "O_O"
$ws.Range("B6").Value = "o_o"
$ws.Range("B7").Value = "bruh"
$ws.Range("B8").Value = 0
$ws.Range("B9").Value = 65
$ws.Range("B10").Value = 65
$ws.Range("B11").Value = 656
$ws.Range("B12").Value = "TL"
$ws.Range("B13").Value = 56
$ws.Range("B14").Value = 56

# --- Insert 4 new rows before the old row 15 ("Kum Döküm") ---
$ws.Rows.Item(15).Resize(4, 1).Insert()

# --- Insert 1 new row before the (now shifted) "Soğuk Maça" row, for "Sıcak Maça" ---
$ws.Rows.Item(22).Insert()

# --- Apply the same label style (bold, bordered, centered) used by column A to all the new/label rows ---
$styleSource = $ws.Range("A2")
$labelRange = $ws.Range("A15:A47")
$styleSource.Copy()
$labelRange.PasteSpecial(-4122)

# --- Fill in labels (column A) and values (column B) for rows 15-47 ---

$ws.Range("A15").Value = "Teslim Türü"
Set-TextValue $ws "B15" "56"
$ws.Range("A16").Value = "Kalıp Göz Sayısı"
$ws.Range("B16").Value = 5656
$ws.Range("A17").Value = "Maça Göz Sayısı"
$ws.Range("B17").Value = 56
$ws.Range("A18").Value = "Kalıp Ömrü Baskı Sayısı"
$ws.Range("B18").Value = 56
$ws.Range("A19").Value = "Kum Döküm"
$ws.Range("B19").Value = 56
$ws.Range("A20").Value = "Kokil Döküm"
$ws.Range("B20").Value = 56
$ws.Range("A21").Value = "Enjeksiyon Döküm"
$ws.Range("B21").Value = 675
$ws.Range("A22").Value = "Sıcak Maça"
$ws.Range("B22").Value = 9
$ws.Range("A23").Value = "Soğuk Maça"
$ws.Range("B23").Value = 798
$ws.Range("A24").Value = "Takalama"
$ws.Range("B24").Value = 657
$ws.Range("A25").Value = "Testere"
$ws.Range("B25").Value = 5
$ws.Range("A26").Value = "Zımpara"
$ws.Range("B26").Value = 57
$ws.Range("A27").Value = "Tesviye"
$ws.Range("B27").Value = 65
$ws.Range("A28").Value = "Kumlama"
$ws.Range("B28").Value = 85
$ws.Range("A29").Value = "Test (Sızdırmazlık)"
$ws.Range("B29").Value = 76
$ws.Range("A30").Value = "Test (Temizleme)"
$ws.Range("B30").Value = 5
$ws.Range("A31").Value = "Isıl İşlem (kg bazında hesaplama için)"
$ws.Range("B31").Value = 5
$ws.Range("A32").Value = "Isıl İşlem (şarf bazında hesaplama için)"
$ws.Range("B32").Value = 765
$ws.Range("A33").Value = "Vibrasyon"
$ws.Range("B33").Value = 0
$ws.Range("A34").Value = "CNC (Dik İşleme)"
$ws.Range("B34").Value = 765
$ws.Range("A35").Value = "CNC Yatay İşleme"
$ws.Range("B35").Value = 765
$ws.Range("A36").Value = "Torna"
$ws.Range("B36").Value = 65
$ws.Range("A37").Value = "Heli-Coil"
$ws.Range("B37").Value = 765
$ws.Range("A38").Value = "Montaj (Parça)"
$ws.Range("B38").Value = 765
$ws.Range("A39").Value = "Montaj (Kaynak)"
$ws.Range("B39").Value = 875
$ws.Range("A40").Value = "Paketleme"
$ws.Range("B40").Value = 8
$ws.Range("A41").Value = "Emprenye"
$ws.Range("B41").Value = 4
$ws.Range("A42").Value = "Kaplama"
$ws.Range("B42").Value = 884
$ws.Range("A43").Value = "Boya"
$ws.Range("B43").Value = 3
$ws.Range("A44").Value = "Dış İşleme"
$ws.Range("B44").Value = 7
$ws.Range("A45").Value = "X-Ray Testi"
$ws.Range("B45").Value = 465
$ws.Range("A46").Value = "Mukavemet Testi"
$ws.Range("B46").Value = 57
$ws.Range("A47").Value = "Nakliye Maliyeti"
$ws.Range("B47").Value = 58
Write-Output "offer_detail sheet updated: details added"
